$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Age" values currently live in C3:C23. Move them one column to the
# left, into B3:B23, leaving C3:C23 empty (B2/C2 headers stay untouched).
$srcRange = $ws.Range("C3:C23")
$dstRange = $ws.Range("B3:B23")

$dstRange.Value2 = $srcRange.Value2
$srcRange.ClearContents()
